# Updated documentation for acknowledgements
#
# Slide 1 ("title slide"): the small credits textbox ("Rectangle 1") gets a
# new centered paragraph acknowledging Eric Callahan (hyperlinked) for
# guidance and support. The textbox auto-sizes (wrap="none" + spAutoFit) to
# fit the new, wider/taller text, so its position/size are adjusted to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "Rectangle 1" credits shape by name (robust to shape ordering).
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($s.Shapes.Item($i).Name -eq "Rectangle 1") {
        $shape = $s.Shapes.Item($i)
        break
    }
}

$tr = $shape.TextFrame.TextRange

# Paragraph 1 already reads "Slide template provided by Sacramento State
# University". Insert a new paragraph right after it (before the existing
# trailing blank paragraph) with the acknowledgement text.
$para1 = $tr.Paragraphs(1)
$null = $para1.InsertAfter("`rAcknowledgement to Eric Callahan for guidance and support")

# Re-fetch the full range and turn the "Eric Callahan" substring into a
# hyperlink run, matching the author/email/colleague acknowledgement pattern
# used elsewhere in this deck.
$linkText = "Eric Callahan"
$full = $shape.TextFrame.TextRange.Text
$linkStart = $full.IndexOf($linkText) + 1
$linkRange = $shape.TextFrame.TextRange.Characters($linkStart, $linkText.Length)
$linkRange.ActionSettings.Item(1).Hyperlink.Address = "https://www.linkedin.com/in/ericcallahan"

# The textbox uses wrap="none" + spAutoFit, so PowerPoint grows/recenters it
# to fit the new (wider and taller) text block. Apply the resulting frame.
$shape.Left = 193.0287401574803
$shape.Width = 338.6753693307087
$shape.Height = 50.89220622440945
